$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 526, shifting existing rows 526-554 down to 527-555
$ws.Rows.Item(526).Insert()

# Populate the newly inserted row 526 with the new price-record data
$ws.Cells.Item(526, 1).Value = 3
$ws.Cells.Item(526, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(526, 3).Value = "Coquimbo"
$ws.Cells.Item(526, 4).Value = 45041
$ws.Cells.Item(526, 5).Value = 5
$ws.Cells.Item(526, 6).Value = 100112040
$ws.Cells.Item(526, 7).Value = "Cilantro"
$ws.Cells.Item(526, 8).Value = "Sin especificar"
$ws.Cells.Item(526, 9).Value = "Primera"
$ws.Cells.Item(526, 10).Value = 185
$ws.Cells.Item(526, 11).Value = 4000
$ws.Cells.Item(526, 12).Value = 4500
$ws.Cells.Item(526, 13).Value = 4243
$ws.Cells.Item(526, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(526, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(526, 16).Value = 1414
$ws.Cells.Item(526, 17).Value = 3
$ws.Cells.Item(526, 18).Value = "Hortaliza"
